# Add unipa transcriptions (-Hasan 026)
#
# The "AF003" event/occurrence transcription series is renamed to the
# "AM003" series:
#   - Sampling Events!A2 (parentEventID): UNIPA-2001ER-AF003        -> UNIPA-2001ER-AM003
#   - Sampling Events!B2 (eventID):       UNIPA-2001ER-AF003-CAPA001 -> UNIPA-2001ER-AM003-CAPA001
#   - Occurrences!A2:A4 (eventID):        UNIPA-2001ER-AF003-CAPA001 -> UNIPA-2001ER-AM003-CAPA001
#   - Occurrences!B2:B4 (occurrenceID):   UNIPA-2001ER-AF003-CAPA001-VE00n -> UNIPA-2001ER-AM003-CAPA001-VE00n
# The "Occurrences" tab also becomes the active tab, with A:B selected on
# both sheets.

$wb = $excel.ActiveWorkbook

$wsSampling = $wb.Worksheets.Item("Sampling Events")
$wsOcc = $wb.Worksheets.Item("Occurrences")

# --- Sampling Events sheet: row 2 ---
$wsSampling.Range("A2").Value = "UNIPA-2001ER-AM003"
$wsSampling.Range("B2").Value = "UNIPA-2001ER-AM003-CAPA001"

# --- Occurrences sheet: rows 2-4 ---
$wsOcc.Range("A2").Value = "UNIPA-2001ER-AM003-CAPA001"
$wsOcc.Range("B2").Value = "UNIPA-2001ER-AM003-CAPA001-VE001"

$wsOcc.Range("A3").Value = "UNIPA-2001ER-AM003-CAPA001"
$wsOcc.Range("B3").Value = "UNIPA-2001ER-AM003-CAPA001-VE002"

$wsOcc.Range("A4").Value = "UNIPA-2001ER-AM003-CAPA001"
$wsOcc.Range("B4").Value = "UNIPA-2001ER-AM003-CAPA001-VE003"

# --- Selections / active sheet ---
$wsSampling.Range("A1:B1048576").Select() | Out-Null
$wsOcc.Activate() | Out-Null
$wsOcc.Range("A1:B1048576").Select() | Out-Null
